$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: Tata Motors (Automobile)
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Automobile"
$ws.Range("C8").Value = "Tata Motors"
$ws.Range("D8").Value = "TATAMOTORS.NS"
$ws.Range("E8").Value = (New-Object System.DateTime(2025, 2, 10))
$ws.Range("F8").Value = "Buy"
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 696.45

# Row 9: Tata Power (Energy)
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Energy"
$ws.Range("C9").Value = "Tata Power"
$ws.Range("D9").Value = "TATAPOWER.NS"
$ws.Range("E9").Value = (New-Object System.DateTime(2025, 2, 10))
$ws.Range("F9").Value = "Buy"
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 359.35

# Match formatting of the existing data rows (date + currency formats, centered)
$ws.Range("E8:E9").NumberFormat = $ws.Range("E7").NumberFormat
$ws.Range("H8:H9").NumberFormat = $ws.Range("H7").NumberFormat
$ws.Range("A8:H9").HorizontalAlignment = -4108
$ws.Range("A8:H9").VerticalAlignment = -4108

# Widen column D slightly to fit the new symbol text, and update selection
$ws.Columns.Item(4).ColumnWidth = 15.44140625
$ws.Range("H10").Select()
